$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 50 and 51: coin name/link order swapped (RocketPoolETH now ranked above FTXToken)
$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("B51").Value = "FTXToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"

# Column D (Price) updates - these are text cells (prices use "." as a thousands
# separator so some are not valid numeric literals). For values that DO look like
# plain numbers, force the Text number format first so Excel keeps them as strings
# (matching the original cell type) instead of silently coercing to a number.
$ws.Range("D2").Value = "36.486.55"
$ws.Range("D3").Value = "1.938.90"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "242.32"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.609"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "56.37"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0806"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.102"
$ws.Range("D12").Value = "2.225.18"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.801"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "13.25"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.12"
$ws.Range("D17").Value = "1.935.66"
$ws.Range("D18").Value = "36.403.84"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "68.86"
$ws.Range("D20").Value = "0.0₃0850"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "226.50"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.92"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.37"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.06"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "159.54"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.131"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.02"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0610"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0982"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.91"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0207"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.13"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "15.70"
$ws.Range("D46").Value = "1.326.93"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "85.45"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.08"
$ws.Range("D50").Value = "2.115.31"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.47"

# Column E (Volume 1h) updates - always text (percent sign + padding spaces
# keep Excel from interpreting them as numbers)
$ws.Range("E2").Value = "  +0.11%  "
$ws.Range("E3").Value = "  -1.01%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("E5").Value = "  -0.86%  "
$ws.Range("E6").Value = "  -1.34%  "
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("E8").Value = "  -3.93%  "
$ws.Range("E9").Value = "  -4.08%  "
$ws.Range("E10").Value = "  -4.47%  "
$ws.Range("E11").Value = "  -1.03%  "
$ws.Range("E12").Value = "  -1.12%  "
$ws.Range("E13").Value = "  -2.93%  "
$ws.Range("E14").Value = "  -4.47%  "
$ws.Range("E15").Value = "  -2.44%  "
$ws.Range("E16").Value = "  -4.44%  "
$ws.Range("E17").Value = "  -2.74%  "
$ws.Range("E18").Value = "  +0.09%  "
$ws.Range("E19").Value = "  -2.04%  "
$ws.Range("E21").Value = "  -2.02%  "
$ws.Range("E22").Value = "  -3.13%  "
$ws.Range("E23").Value = "  -0.05%  "
$ws.Range("E24").Value = "  -6.16%  "
$ws.Range("E25").Value = "  -0.05%  "
$ws.Range("E26").Value = "  -5.05%  "
$ws.Range("E27").Value = "  -3.10%  "
$ws.Range("E28").Value = "  +7.42%  "
$ws.Range("E29").Value = "  -3.34%  "
$ws.Range("E30").Value = "  -1.10%  "
$ws.Range("E31").Value = "  -7.29%  "
$ws.Range("E32").Value = "  -4.19%  "
$ws.Range("E33").Value = "  -4.74%  "
$ws.Range("E34").Value = "  -5.38%  "
$ws.Range("E35").Value = "  -0.12%  "
$ws.Range("E36").Value = "  -1.24%  "
$ws.Range("E37").Value = "  -1.31%  "
$ws.Range("E38").Value = "  +0.76%  "
$ws.Range("E39").Value = "  +9.85%  "
$ws.Range("E40").Value = "  -0.11%  "
$ws.Range("E41").Value = "  +1.11%  "
$ws.Range("E42").Value = "  -1.72%  "
$ws.Range("E43").Value = "  -4.68%  "
$ws.Range("E44").Value = "  -0.75%  "
$ws.Range("E45").Value = "  -2.99%  "
$ws.Range("E46").Value = "  -1.52%  "
$ws.Range("E47").Value = "  -4.14%  "
$ws.Range("E48").Value = "  -5.10%  "
$ws.Range("E49").Value = "  -0.56%  "
$ws.Range("E50").Value = "  -1.22%  "
$ws.Range("E51").Value = "  +14.45%  "

